# Consolidate the "An" / " " / "image" runs in the caption textbox into a
# single run containing "An image".
#
# Setting TextRange.Text to the value it already holds is a no-op (the
# runtime only rewrites runs when the text actually changes), so first set
# it to a placeholder value to force the rewrite, then set it to the final
# consolidated text. This produces a single run instead of reintroducing
# per-character/lang-stamped runs.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 3")
$sh.TextFrame.TextRange.Text = "x"
$sh.TextFrame.TextRange.Text = "An image"
